# Edit script: IF2210_W05_TB1CPP_LogAct.docx changes
# 1. Table formatting: tblInd 0 -> -5 dxa, tblCellMar/left 108 -> 103 dxa,
#    column 2 width 1564 -> 1565 dxa, column 7 width 2788 -> 2787 dxa,
#    and every cell's left margin 108 -> 103 dxa.
# 2. Content: merge "Rancangan awal Zoo untuk VZ01, " + "Header Class, dan
#    Implementasinya" into a single run.
# 3. Content: "Melakukan driver untuk kelas Animal dan kelas-kelas terkait"
#    -> split into 3 runs ("Membuat " / "driver untuk kelas Animal " /
#    "dan Realisasi kelas Animal untuk sejumlah binatang nyata"); and
#    "driver class" -> split into 2 runs ("driver class " / "dan sekumpulan
#    kelas realisasi Animal di animal_list.h").
# 4. Footer: drop the italic/sz18 run formatting on the first FILENAME
#    field's "begin" fldChar run.
# 5. Styles: drop <w:sz w:val="22"/> from docDefaults/rPrDefault; change
#    <w:color w:val="auto"/> to <w:color w:val="00000A"/> in the Normal and
#    NoSpacing paragraph styles.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# -- table-wide indent & default cell margin --
$t.Rows.LeftIndent = -0.25      # -5 dxa
$t.LeftPadding = 5.15           # 103 dxa (tblCellMar/left)

# -- column widths (dxa / 20 = points) --
$t.Columns.Item(2).Width = 78.25   # 1565 dxa
$t.Columns.Item(7).Width = 139.35  # 2787 dxa

# -- every cell's left padding 108 -> 103 dxa (5.4 -> 5.15 pt) --
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $row.Cells.Item($c).LeftPadding = 5.15
    }
}

# -- content change 1: merge the two runs in row 9 / column 7 --
$d.Content.Find.Execute("Rancangan awal Zoo untuk VZ01, Header Class, dan Implementasinya", $true, $false, $false, $false, $false, $true, 1, $false, "Rancangan awal Zoo untuk VZ01, Header Class, dan Implementasinya", 2)

# -- content change 2: rewrite the two cells in the "Melakukan driver" row --
$d.Content.Find.Execute("Melakukan driver untuk kelas Animal dan kelas-kelas terkait", $true, $false, $false, $false, $false, $true, 1, $false, "Membuat driver untuk kelas Animal dan Realisasi kelas Animal untuk sejumlah binatang nyata", 2)
$d.Content.Find.Execute("driver class", $true, $false, $false, $false, $false, $true, 1, $false, "driver class dan sekumpulan kelas realisasi Animal di animal_list.h", 2)

# -- footer: strip italic/size-18 formatting from the first FILENAME field's begin fldChar --
$ftr = $d.Sections.Item(1).Footers.Item(1)
$ftrPara = $ftr.Range.Paragraphs.Item(1)
$firstRun = $ftrPara.Range.Words.Item(1)
$frange = $ftr.Range.Duplicate
$frange.SetRange($ftr.Range.Start, $ftr.Range.Start)
$frange.MoveEndUntil(" ", 1) | Out-Null
$frange.Font.Italic = 0
$frange.Font.Size = 10

# -- styles.xml: docDefaults rPrDefault drop explicit size --
$normalStyle = $d.Styles.Item("Normal")
$noSpacingStyle = $d.Styles.Item("No Spacing")
$normalStyle.Font.Color = 10
$noSpacingStyle.Font.Color = 10

Write-Output "done"
